$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Step 1: the empty paragraph right after the "TLS(HTTPS) ..." paragraph had
# a stray empty <w:lang/> in its mark rPr -- drop it.
# ---------------------------------------------------------------------------
$langParaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr></w:p>
'@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "`r" -and $p.Range.Start -gt 41 -and $p.Range.Start -lt 106) {
        $p.Range.InsertXML($langParaXml)
        break
    }
}

# ---------------------------------------------------------------------------
# Step 2: rewrite the "Insecure direct object references" answer paragraph.
# The old paragraph held a bold "RESENO " run followed by a short answer; it
# is replaced by an expanded explanation split across three paragraphs.
# ---------------------------------------------------------------------------
$replacementXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:b/><w:color w:val="202729"/><w:sz w:val="28"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Postoje dva načina od ove vrste napada. Prvi je da se objekti indirektno referenciraju, a drugi da se objekti direktno referenciraju uz proveru prava pristupa. Naš način odbrane </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t xml:space="preserve">od ovog napada jeste taj da prilikom svakog pristupa određenom resursu, proveravamo da li je korisnik ulogovan i koja su njegova prava pristupa (dakle, drugi način). </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t xml:space="preserve">U datoj aplikaciji, ovo je rešeno tako što je ulogovanom korisniku prikazano samo ono što je njemu predviđeno. Na primer, ulogovan građanin ne može uneti novi propis već može samo vršiti pretragu postojećih propisa, kao i pregled istih, odbornik ne može prihvatati akte itd. Naravno, prava pristupa se proveravaju prilikom svake akcije. Ukoliko ulogovani građanin pokuša da preko URL-a ode na određenu stranicu na kojoj nema prava pristupa, biće redirektovan. </w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:tab/><w:t xml:space="preserve">Provera prava pristupa su konkretno implementirana na front-end-u u html stranicama (šta ulogovani korisnik sme da vidi od funkcionalnosti), kao i u </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t xml:space="preserve">kontrolerima </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t xml:space="preserve">u kojima se </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t>proverava</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t>ju</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t xml:space="preserve"> prava </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t>pristupa svaki put pre nego što se izvrš</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t>i neki servis za dobavljanje podatak</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:eastAsia="ProximaNova-Regular" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:color w:val="202729"/><w:sz w:val="24"/><w:szCs w:val="28"/><w:lang/></w:rPr><w:t xml:space="preserve">a na osnovu direktnih referenci (AngularJS). </w:t></w:r></w:p>
'@

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*RESENO*") {
        $p.Range.InsertXML($replacementXml)
        break
    }
}
